$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as text (matching the
# original inlineStr cells) instead of auto-converting numeric-looking
# strings (e.g. "6.10", "0.670") into floating point numbers.

# --- Per-row Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "'63.907.81"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "'3.075.52"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D5").Value = "'558.37"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'142.99"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'3.074.06"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "  -6.71%  "
$ws.Range("D12").Value = "'0.469"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'34.92"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "'3.583.98"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'63.903.83"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "'3.075.73"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "'480.03"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("D24").Value = "'14.04"
$ws.Range("E24").Value = "  +11.12%  "
$ws.Range("D25").Value = "'81.08"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "'7.93"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "'2.06"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'26.26"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").Value = "'2.47"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").Value = "'5.56"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "'55.77"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'6.17"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'448.59"
$ws.Range("E37").Value = "  -4.07%  "
$ws.Range("D40").Value = "'0.0816"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'2.965.80"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").Value = "'8.21"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  -4.76%  "
$ws.Range("D44").Value = "'27.78"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "'0.260"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("D48").Value = "'0.112"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "'119.23"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'0.0₃0510"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("E51").Value = "  -0.78%  "

# --- Rows 38/39: dogwifhat jumps above VeChain in the rankings ---
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'2.95"
$ws.Range("E38").Value = "  +15.16%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0405"
$ws.Range("E39").Value = "  +1.94%  "
